$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Name value (B4) moves into Title value (B5); Name value cell is cleared.
$nameValue = $meta.Range("B4").Value()
$meta.Range("B5").Value = $nameValue
$meta.Range("B4").ClearContents()

# Date value (B8) updated
$meta.Range("B8").Value = "2026-01-07T15:20:53+00:00"

# --- Mapping Table 0 sheet ---
$map0 = $wb.Worksheets.Item("Mapping Table 0")
for ($r = 9; $r -le 13; $r++) {
    $cell = $map0.Cells.Item($r, 4)
    $old = $cell.Value()
    $new = $old.Replace("entryRelationship.", "entryRelationship:")
    $cell.Value = $new
}

# --- Mapping Table 1 sheet ---
$map1 = $wb.Worksheets.Item("Mapping Table 1")
for ($r = 9; $r -le 13; $r++) {
    $cell = $map1.Cells.Item($r, 1)
    $old = $cell.Value()
    $new = $old.Replace("entryRelationship.", "entryRelationship:")
    $cell.Value = $new
}
